# "Trying in OOP - Dict window"
# Adds a new "Sheet1" (a simple pop-up style dictionary window) between the
# "word" and "multiply" sheets, adds a batch of new food vocabulary rows to
# the "word" sheet (rows 84-93, previously blank placeholder rows), fills in
# the missing "it/this" <-> "oно" translation on row 14, and drops the stale
# AutoFilter sort-state left over from an old sort.

$wb = $excel.ActiveWorkbook
$wordSheet = $wb.Worksheets.Item("word")
$multiplySheet = $wb.Worksheets.Item("multiply")

# ---------------------------------------------------------------------
# 1. New vocabulary: fill in previously-empty placeholder rows 84-93 with
#    food words (English in B, Russian in D, an extra English synonym in C
#    for row 86), then unhide those rows since they now carry live data.
# ---------------------------------------------------------------------
$wordSheet.Range("A84:A93").EntireRow.Hidden = $false

$foodRows = 84..93
$englishWords = @("meat","potato","biscuits","cake","orange juice","pasta","carrots","sausages","rice","popcorn")
$russianWords = @("мясо","картошка","печенье","торт","апельсиновый сок","макароны","морковка","сосиски","рис","попкорн")

for ($i = 0; $i -lt $foodRows.Length; $i++) {
    $r = $foodRows[$i]
    $wordSheet.Cells.Item($r, 1).Value = 1
    $wordSheet.Cells.Item($r, 2).Value = $englishWords[$i]
}
for ($i = 0; $i -lt $foodRows.Length; $i++) {
    $r = $foodRows[$i]
    $wordSheet.Cells.Item($r, 4).Value = $russianWords[$i]
    if ($r -eq 86) {
        $wordSheet.Cells.Item($r, 3).Value = "cookie"
    }
}

# row 14 ("it" / "это") was missing the synonym "оно" in column E
$wordSheet.Cells.Item(14, 5).Value = "оно"

# ---------------------------------------------------------------------
# 2. Drop the stale AutoFilter sort-state (column-1 filter stays as-is).
# ---------------------------------------------------------------------
$wordSheet.AutoFilterMode = $false
$wordSheet.Range("A1:F274").AutoFilter(1, @("1"))

# ---------------------------------------------------------------------
# 3. New "Sheet1" - a small popup/dictionary window: column A lists the
#    Russian word to translate (big 14pt font), column B is a blank
#    answer box with a bottom border rule.
# ---------------------------------------------------------------------
$dict = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wordSheet)
$dict.Name = "Sheet1"

$dictWords = @("я","вы","ваш","он","его","она","её","это","мы","наш","они","их","мой",
               "мясо","картошка","печенье","торт","апельсиновый сок","макароны","морковка","сосиски","рис","попкорн")

for ($i = 0; $i -lt $dictWords.Length; $i++) {
    $row = $i + 1
    $cellA = $dict.Cells.Item($row, 1)
    $cellA.Value = $dictWords[$i]
    $cellA.Font.Size = 14
    $dict.Rows.Item($row).RowHeight = 18.75

    $cellB = $dict.Cells.Item($row, 2)
    $cellB.Borders.Item(9).LineStyle = 1
}

$dict.Columns.Item(1).ColumnWidth = 22.86
$dict.Columns.Item(2).ColumnWidth = 28.43

# ---------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping, matching the saved UI state.
# ---------------------------------------------------------------------
$wordSheet.Activate()
$wordSheet.Range("D7:D93").Select()

$dict.Activate()
$dict.Range("E8").Select()

Write-Host "done"
